$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.200.11"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.852.71"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "0.6989"
$ws.Range("E5").Value = "  +1.42%  "

$ws.Range("D6").Value = "237.72"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "0.07886"
$ws.Range("E8").Value = "  +1.43%  "

$ws.Range("D9").Value = "0.3017"
$ws.Range("E9").Value = "  -1.05%  "

$ws.Range("D10").Value = "23.72"
$ws.Range("E10").Value = "  +2.21%  "

$ws.Range("D11").Value = "0.08111"
$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("D12").Value = "1.855.15"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").Value = "5.176"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("D14").Value = "0.7045"
$ws.Range("E14").Value = "  -2.30%  "

$ws.Range("D15").Value = "89.51"
$ws.Range("E15").Value = "  +0.23%  "

$ws.Range("D16").Value = "29.217.38"
$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").Value = "5.802"
$ws.Range("E17").Value = "  +1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007818"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.20"
$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").Value = "235.72"
$ws.Range("E20").Value = "  +0.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").Value = "2.100.85"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").Value = "7.501"
$ws.Range("E24").Value = "  +0.53%  "

$ws.Range("D25").Value = "162.42"
$ws.Range("E25").Value = "  +0.22%  "

$ws.Range("D26").Value = "8.861"
$ws.Range("E26").Value = "  -1.18%  "

$ws.Range("D27").Value = "0.1426"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").Value = "18.02"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("E29").Value = "  -1.87%  "

$ws.Range("D30").Value = "1.405"
$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("D31").Value = "1.474"
$ws.Range("E31").Value = "  -0.55%  "

$ws.Range("D32").Value = "4.324"
$ws.Range("E32").Value = "  -4.21%  "

$ws.Range("D33").Value = "4.008"
$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").Value = "0.05155"
$ws.Range("E34").Value = "  -0.97%  "

$ws.Range("D35").Value = "1.163"
$ws.Range("E35").Value = "  -1.64%  "

$ws.Range("D36").Value = "0.7105"
$ws.Range("E36").Value = "  +1.10%  "

$ws.Range("D37").Value = "0.9984"
$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("D38").Value = "2.678"
$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("D39").Value = "0.01845"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  +0.98%  "

$ws.Range("D41").Value = "1.149.92"
$ws.Range("E41").Value = "  +5.34%  "

$ws.Range("D42").Value = "0.9243"
$ws.Range("E42").Value = "  -1.79%  "

$ws.Range("D43").Value = "5.964"
$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").Value = "0.4235"
$ws.Range("E44").Value = "  -1.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.10"
$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").Value = "102.93"
$ws.Range("E47").Value = "  +0.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5290"
$ws.Range("E48").Value = "  -2.92%  "

$ws.Range("D49").Value = "1.735"
$ws.Range("E49").Value = "  -3.33%  "

$ws.Range("D50").Value = "9.134"
$ws.Range("E50").Value = "  -0.30%  "

$ws.Range("D51").Value = "6.952"
$ws.Range("E51").Value = "  -0.76%  "
